# Insert a new row at position 20, shifting the existing rows 20-24 down to 21-25,
# then populate the new row 20 with the latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(20).Insert()

$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44642
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100104
$ws.Cells.Item(20, 8).Value = "Frutos de pepita"
$ws.Cells.Item(20, 9).Value = 100104005
$ws.Cells.Item(20, 10).Value = "Pera"
$ws.Cells.Item(20, 11).Value = "Packham's Triumph"
$ws.Cells.Item(20, 12).Value = "Segunda"
$ws.Cells.Item(20, 13).Value = 270
$ws.Cells.Item(20, 14).Value = 19000
$ws.Cells.Item(20, 15).Value = 20000
$ws.Cells.Item(20, 16).Value = 19500
$ws.Cells.Item(20, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 1083
$ws.Cells.Item(20, 20).Value = 18
